$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row (new row 9): "Temps de mise en oeuvre..." risk ---
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).ClearFormats()

$ws.Range("F9").Value = "Revoir le cahier de charges"
$ws.Range("B9").Value = "Temps de mise en œuvre prévue d'une fonctionnalité logicielle trop élevée"
$ws.Range("G9").Value = "Douglas R."
$ws.Range("G9").WrapText = $true
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 15
$ws.Range("H9").Value = "2j"

# --- Insert second new row (new row 14): "Problèmes inattendus à débogger" risk ---
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).ClearFormats()

$ws.Range("B14").Value = "Problèmes inattendus à débogger"
$ws.Range("F14").Value = "Revoir le Gantt et allouer plus de temps au déboggage"
$ws.Range("G14").Value = "Badr B." + [char]10 + "Douglas R."
$ws.Range("G14").WrapText = $true
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 5
$ws.Range("H14").Value = "1j"
$ws.Rows.Item(14).RowHeight = 30

# --- Update selection / view state to match final state ---
$ws.Range("F22").Select()

Write-Output "done"
